# Auto-generated script to update cached market-price values
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 998.3333
$ws.Range("I9").Value = 998.3333
$ws.Range("K9").Value = 998.3333
$ws.Range("M9").Value = -829.3333
$ws.Range("H33").Value = 466
$ws.Range("I33").Value = 466
$ws.Range("K33").Value = 466
$ws.Range("M33").Value = -237
$ws.Range("H40").Value = 1342.5714
$ws.Range("I40").Value = 1255.1111
$ws.Range("K40").Value = 1255.1111
$ws.Range("M40").Value = -1080.1111
$ws.Range("H64").Value = 4749.5
$ws.Range("J64").Value = 4499.5
$ws.Range("L64").Value = 4499.5
$ws.Range("N64").Value = -4995.5
$ws.Range("H67").Value = 4749.5
$ws.Range("J67").Value = 4499.5
$ws.Range("L67").Value = 4499.5
$ws.Range("N67").Value = -6215.5
$ws.Range("H70").Value = 4361.8887
$ws.Range("I70").Value = 3175.5
$ws.Range("J70").Value = 5116.864
$ws.Range("K70").Value = 9526.5
$ws.Range("L70").Value = 15350.592
$ws.Range("M70").Value = -9256.5
$ws.Range("N70").Value = -15890.592
$ws.Range("H73").Value = 4361.8887
$ws.Range("I73").Value = 3175.5
$ws.Range("J73").Value = 5116.864
$ws.Range("K73").Value = 9526.5
$ws.Range("L73").Value = 15350.592
$ws.Range("M73").Value = -8590.5
$ws.Range("N73").Value = -17222.592
$ws.Range("H132").Value = 5865.5835
$ws.Range("I132").Value = 6262.5454
$ws.Range("K132").Value = 18787.6362
$ws.Range("M132").Value = -16257.6362
$ws.Range("H135").Value = 698.8125
$ws.Range("I135").Value = 698.8125
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6289.3125
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3754.3125
$ws.Range("H137").Value = 3309.4285
$ws.Range("J137").Value = 3122.6667
$ws.Range("L137").Value = 9368.000100000001
$ws.Range("N137").Value = -14468.0001
$ws.Range("N135").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10460.8
$ws.Range("I32").Value = 10480.042
$ws.Range("K32").Value = 10480.042
$ws.Range("M32").Value = -10193.042
$ws.Range("H61").Value = 2970
$ws.Range("I61").Value = 2767.6155
$ws.Range("K61").Value = 2767.6155
$ws.Range("M61").Value = -2555.6155
$ws.Range("H74").Value = 2028.5
$ws.Range("I74").Value = 1997.2
$ws.Range("K74").Value = 1997.2
$ws.Range("M74").Value = -1123.2
$ws.Range("H77").Value = 2028.5
$ws.Range("I77").Value = 1997.2
$ws.Range("K77").Value = 9986
$ws.Range("M77").Value = -5618
$ws.Range("H110").Value = 5788.8
$ws.Range("I110").Value = 5702.2856
$ws.Range("K110").Value = 5702.2856
$ws.Range("M110").Value = -3657.2856
$ws.Range("H122").Value = 1837.1052
$ws.Range("I122").Value = 1911.2778
$ws.Range("J122").Value = 502
$ws.Range("K122").Value = 5733.8334
$ws.Range("L122").Value = 1506
$ws.Range("M122").Value = -3283.8334
$ws.Range("N122").Value = -6406
$ws.Range("H132").Value = 1734.0555
$ws.Range("I132").Value = 1734.0555
$ws.Range("K132").Value = 5202.166499999999
$ws.Range("M132").Value = -2672.166499999999
$ws.Range("H136").Value = 2970
$ws.Range("I136").Value = 2767.6155
$ws.Range("K136").Value = 8302.8465
$ws.Range("M136").Value = -5752.8465

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 834.8182
$ws.Range("I107").Value = 773.25
$ws.Range("J107").Value = 999
$ws.Range("K107").Value = 773.25
$ws.Range("L107").Value = 999
$ws.Range("M107").Value = 1146.75
$ws.Range("N107").Value = -4839

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 913445.75
$ws.Range("I19").Value = 1428899.6
$ws.Range("J19").Value = 11401.5
$ws.Range("K19").Value = 1428899.6
$ws.Range("L19").Value = 11401.5
$ws.Range("M19").Value = -1428729.6
$ws.Range("N19").Value = -11741.5
$ws.Range("H24").Value = 913445.75
$ws.Range("I24").Value = 1428899.6
$ws.Range("J24").Value = 11401.5
$ws.Range("K24").Value = 1428899.6
$ws.Range("L24").Value = 11401.5
$ws.Range("M24").Value = -1428729.6
$ws.Range("N24").Value = -11741.5
$ws.Range("H88").Value = 12080.667
$ws.Range("J88").Value = 12080.667
$ws.Range("L88").Value = 12080.667
$ws.Range("N88").Value = -12892.667
$ws.Range("H91").Value = 12080.667
$ws.Range("J91").Value = 12080.667
$ws.Range("L91").Value = 12080.667
$ws.Range("N91").Value = -14888.667
$ws.Range("H134").Value = 3803.7273
$ws.Range("I134").Value = 3684.1
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 11052.3
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -8517.299999999999
$ws.Range("N134").Value = -20070

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1928621.8
$ws.Range("I4").Value = 57196.465
$ws.Range("J4").Value = 30000000
$ws.Range("K4").Value = 171589.395
$ws.Range("L4").Value = 90000000
$ws.Range("M4").Value = -171477.395
$ws.Range("N4").Value = -90000224
$ws.Range("H10").Value = 164.36363
$ws.Range("I10").Value = 164.36363
$ws.Range("K10").Value = 493.09089
$ws.Range("M10").Value = -354.09089
$ws.Range("H12").Value = 337.6154
$ws.Range("I12").Value = 310.22223
$ws.Range("K12").Value = 930.66669
$ws.Range("M12").Value = -757.66669
$ws.Range("H14").Value = 388.91666
$ws.Range("I14").Value = 388.91666
$ws.Range("K14").Value = 1166.74998
$ws.Range("M14").Value = -993.7499800000001
$ws.Range("H109").Value = 5343.909
$ws.Range("I109").Value = 662
$ws.Range("K109").Value = 1986
$ws.Range("M109").Value = -946

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("J5").Value = 10000
$ws.Range("L5").Value = 10000
$ws.Range("N5").Value = -10224
$ws.Range("H122").Value = 1650.875
$ws.Range("I122").Value = 1650.875
$ws.Range("K122").Value = 4952.625
$ws.Range("M122").Value = -2502.625
$ws.Range("H132").Value = 3460
$ws.Range("I132").Value = 2950.8
$ws.Range("J132").Value = 6006
$ws.Range("K132").Value = 8852.400000000001
$ws.Range("L132").Value = 18018
$ws.Range("M132").Value = -6322.400000000001
$ws.Range("N132").Value = -23078

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 21458.572
$ws.Range("I2").Value = 13105.5
$ws.Range("J2").Value = 24799.8
$ws.Range("K2").Value = 13105.5
$ws.Range("L2").Value = 24799.8
$ws.Range("M2").Value = -12993.5
$ws.Range("N2").Value = -25023.8
$ws.Range("H7").Value = 9999.5
$ws.Range("I7").Value = 9999.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 9999.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -9887.5
$ws.Range("H40").Value = 2474.5
$ws.Range("I40").Value = 2474.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2474.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2338.5
$ws.Range("H55").Value = 781.4286
$ws.Range("J55").Value = 1047.8889
$ws.Range("L55").Value = 1047.8889
$ws.Range("N55").Value = -1393.8889
$ws.Range("H68").Value = 3759
$ws.Range("I68").Value = 2265.3333
$ws.Range("K68").Value = 2265.3333
$ws.Range("M68").Value = -1516.3333
$ws.Range("H71").Value = 3759
$ws.Range("I71").Value = 2265.3333
$ws.Range("K71").Value = 11326.6665
$ws.Range("M71").Value = -7582.666499999999
$ws.Range("H122").Value = 7500.8
$ws.Range("I122").Value = 7001.3335
$ws.Range("K122").Value = 21004.0005
$ws.Range("M122").Value = -18554.0005
$ws.Range("H126").Value = 9999.5
$ws.Range("I126").Value = 9999.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 29998.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -27528.5
$ws.Range("N7").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("H100").Value = 1029
$ws.Range("I100").Value = 286.25
$ws.Range("K100").Value = 572.5
$ws.Range("M100").Value = -31.5
$ws.Range("H103").Value = 24475
$ws.Range("J103").Value = 24475
$ws.Range("L103").Value = 24475
$ws.Range("N103").Value = -26819
$ws.Range("H122").Value = 3287
$ws.Range("I122").Value = 2922.7917
$ws.Range("K122").Value = 8768.375100000001
$ws.Range("M122").Value = -6318.375100000001
$ws.Range("H126").Value = 3715.8333
$ws.Range("I126").Value = 3715.8333
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 11147.4999
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8677.499899999999
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("N126").ClearContents()
